$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B. This shifts the existing
# B,C,D,E (Jun_17, Jun_15, Jun_13, Jun_10) report columns right to
# E,F,G,H respectively, preserving their per-cell formatting (e.g. the
# yellow "latest rating" highlight fill moves along with the data).
$ws.Range("B1:B27").EntireColumn.Insert()
$ws.Range("B1:B27").EntireColumn.Insert()
$ws.Range("B1:B27").EntireColumn.Insert()

# Give the three new columns the same custom width as the others.
$ws.Range("B1:D29").ColumnWidth = 7.1667

# New date-group headers for the newly inserted columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Default filler value ("UN") for the new columns' data rows, matching
# the convention used throughout the rest of the sheet.
$ws.Range("B2:D21").Value = "UN"
$ws.Range("B22:B22").Value = "UN"
$ws.Range("B23:D27").Value = "UN"
$ws.Range("E22").Value = "UN"

# Row 22 (BidaskClub) gets the new 6/22/2018 upgrade note in both of
# the new Jun_26 columns.
$ws.Range("C22").Value = "6/22/2018,Upgrades,Hold -> Buy,"
$ws.Range("D22").Value = "6/22/2018,Upgrades,Hold -> Buy,"

# Row 27 (Barclays) gets the new 6/27/2018 initiation note.
$ws.Range("B27").Value = "6/27/2018,Initiates,Overweight,`$74.00"

# New analyst-group rows appended at the bottom.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
